$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 17) with the latest trade data, mirroring the
# existing rows' structure (A:BuyPrice-history, B:StartPrinciple, C:BuyPrice,
# D:SellPrice, E:IsShortSell, F:Price Change %, G:Date, H:Profitable)
$ws.Range("A17").Value = 9516.94
$ws.Range("B17").Value = 9801.17
$ws.Range("C17").Value = 294.14
$ws.Range("D17").Value = 302.66000000000003
$ws.Range("E17").Value = $true
$ws.Range("F17").Value = 2.9
$ws.Range("G16").Copy() | Out-Null
$ws.Range("G17").PasteSpecial(-4122) | Out-Null
$ws.Range("G17").Value = 42626.544363425928
$ws.Range("H17").Value = $false
